$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 16, shifting existing rows 16.. down by one
# (mirrors -1 dimension change from A1:R123 to A1:R124)
$ws.Rows("16:16").Insert()

# Copy the (now shifted-down) original row 16 content into the new row 16,
# since the new row was a duplicate entry of that row with a few fields updated.
$src = 17
$dst = 16

$ws.Cells.Item($dst, 1).Value2 = $ws.Cells.Item($src, 1).Value2    # A Mercado ID
$ws.Cells.Item($dst, 2).Value2 = $ws.Cells.Item($src, 2).Value2    # B Mercado
$ws.Cells.Item($dst, 3).Value2 = $ws.Cells.Item($src, 3).Value2    # C Region
$ws.Cells.Item($dst, 4).Value2 = 44473                              # D Fecha (new)
$ws.Cells.Item($dst, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat
$ws.Cells.Item($dst, 5).Value2 = $ws.Cells.Item($src, 5).Value2    # E Codreg
$ws.Cells.Item($dst, 6).Value2 = $ws.Cells.Item($src, 6).Value2    # F Categoria ID
$ws.Cells.Item($dst, 7).Value2 = $ws.Cells.Item($src, 7).Value2    # G Categoria
$ws.Cells.Item($dst, 8).Value2 = $ws.Cells.Item($src, 8).Value2    # H Variedad
$ws.Cells.Item($dst, 9).Value2 = $ws.Cells.Item($src, 9).Value2    # I Calidad
$ws.Cells.Item($dst, 10).Value2 = 400                               # J Volumen (new)
$ws.Cells.Item($dst, 11).Value2 = 40000                             # K Precio minimo (new)
$ws.Cells.Item($dst, 12).Value2 = 41000                             # L Precio maximo (new)
$ws.Cells.Item($dst, 13).Value2 = 40500                             # M Precio promedio ponderado (new)
$ws.Cells.Item($dst, 14).Value2 = $ws.Cells.Item($src, 14).Value2  # N Unidad de comercializacion
$ws.Cells.Item($dst, 15).Value2 = $ws.Cells.Item($src, 15).Value2  # O Origen
$ws.Cells.Item($dst, 16).Value2 = 1620                              # P Precio $/Kg (new)
$ws.Cells.Item($dst, 17).Value2 = $ws.Cells.Item($src, 17).Value2  # Q Kg o Unidades
$ws.Cells.Item($dst, 18).Value2 = $ws.Cells.Item($src, 18).Value2  # R Clasificacion
